# "paragrafo 1 e 2"
#
# Before:
#   P1: "Maus-tratos aos animais" (title)
#   P2: empty paragraph that only carries the _GoBack bookmark
#
# After:
#   P1: "Maus-tratos aos animais" (title)                       -- unchanged
#   P2 (new): "Talvez por falta de informacao ... nas ruas do Brasil."
#   P3 (was P2): "Quando se fala de maus-tratos ... mesmo com o e" +
#                the (untouched, still in place) _GoBack bookmark +
#                "sforco de entidades assistenciais ... deve ser reconhecido."

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Step 1: split off a brand-new, empty paragraph right before the
# paragraph that holds the _GoBack bookmark. That new paragraph
# becomes the document's paragraph 2.
# ------------------------------------------------------------------
$bookmarkPara = $d.Paragraphs(2)
$bookmarkPara.Range.InsertParagraphBefore()

# ------------------------------------------------------------------
# Step 2: populate the new paragraph 2 with the "Talvez por falta..."
# text, built up from the same four sentence chunks the original
# author typed, kept as four distinct runs (like in the target
# document) by temporarily parking a bookmark at the tail end of the
# paragraph while we prepend each chunk with InsertBefore -- the
# engine only coalesces a freshly inserted run into its neighbour
# when nothing else anchors a run boundary between them, so the
# scratch bookmark keeps the chunks from collapsing into one run.
# The scratch bookmark itself is removed once the text is in place.
# ------------------------------------------------------------------
$para2 = $d.Paragraphs(2)
$para2.Range.Text = " máximo dele? Há espaço em casa para ter um animal de estimação? Tenho capacidade de oferecer uma qualidade de vida boa para ele? Essas são algumas perguntas que muitas vezes não são consideradas antes da aquisição do animal. Seja um gato, um cachorro, temos que levar em consideração uma série de pontos para que o animal não seja uma futura vítima de maus-tratos nas ruas do Brasil."

$anchorPos = $d.Paragraphs(2).Range.End
$anchorRange = $d.Range($anchorPos - 1, $anchorPos - 1)
$d.Bookmarks.Add("ScratchAnchor2", $anchorRange)

$d.Paragraphs(2).Range.InsertBefore("Qual será o tamanho")
$d.Paragraphs(2).Range.InsertBefore(" ou comprados sem antes terem sido analisados pelos compradores. O que será gasto com o animal? ")
$d.Paragraphs(2).Range.InsertBefore("Talvez por falta de informação de muitas pessoas, muitos animais domésticos são adquiridos")

$d.Bookmarks("ScratchAnchor2").Delete()

# ------------------------------------------------------------------
# Step 3: the paragraph that still carries the _GoBack bookmark is now
# paragraph 3. Add the trailing ("...sforço...") text AFTER the
# bookmark first -- while the paragraph is still otherwise empty this
# reliably lands past bookmarkEnd -- and only then add the leading
# ("Quando se fala...") text BEFORE the bookmark, so bookmarkStart /
# bookmarkEnd stay exactly where they were, splitting "esforço" into
# "e" + "sforço" around them, just like in the target document.
# ------------------------------------------------------------------
$para3Range = $d.Paragraphs(3).Range
$para3Range.InsertAfter("sforço de entidades assistenciais e colaboradores, esforço esse que deve ser reconhecido.")

$para3Range = $d.Paragraphs(3).Range
$para3Range.InsertBefore("remover os animais das ruas é tanto para proteger os animais, quanto para proteger os seres humanos. Animais de rua possuem um grande risco de contrair uma doença, já que o animal pelo seu instinto, independente da raça, tem a tendência de revirar lixos e comer comida jogada na rua. Por esse motivo, muitos dos animais chegam ao CCZ com alguma doença grave, sendo necessário na maioria dos casos encaminhar o animal para a eutanásia (ato de proporcionar morte sem sofrimento a pacientes em estado terminal). É minoria o número de animais que saem vivos e com abrigo de um CCZ, mesmo com o e")

$para3Range = $d.Paragraphs(3).Range
$para3Range.InsertBefore("Quando se fala de maus-tratos, se inclui também o abandono do animal. Animais abandonados vão para o CCZ – Centro de Controle de Zoonoses. Zoonoses são doenças que podem ser transmitidas de um animal para um ser humano, ou de um ser humano para um animal. Portanto, o ato de ")

Write-Host "Paragraph count:" $d.Paragraphs.Count
Write-Host "Bookmarks left:" $d.Bookmarks.Count
